$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    4   = -7.093999999999998
    6   = -7.766699999999999
    7   = -7.817400000000001
    16  = -8.625700000000002
    20  = -7.422999999999999
    28  = -8.333899999999996
    29  = -7.252299999999999
    32  = -9.153399999999994
    40  = -8.192799999999995
    46  = -8.195799999999997
    51  = -7.632799999999999
    52  = -7.576199999999997
    57  = -8.199499999999999
    59  = -8.322299999999997
    62  = -9.148699999999989
    66  = -7.212800000000003
    73  = -7.574999999999996
    74  = -8.233400000000001
    92  = -6.189500000000002
    100 = -8.287499999999998
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
